# Generate Report for Handback
# Updates the timestamp values on the handback-status report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-08-26 11:03:38"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn.Range("H2").Value = "2016-08-26 11:03:33"
$wsZhCn.Range("K2").Value = "2016-08-26 11:03:51"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsDeDe.Range("H2").Value = "2016-08-26 11:03:38"
$wsDeDe.Range("K2").Value = "2016-08-26 11:03:58"
